$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K -> new F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting/styles from the (now-shifted) original D:E columns
# (now at F:G) onto the freshly inserted D:E columns so the new cells carry
# the same style index as their neighbours (reuses existing cellXfs entries).
# Done per contiguous data block so blank separator rows (36/37, 78/79) are
# left untouched.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarterly columns (D, E) with the new figures
$ws.Cells.Item(7, 4).Value2 = 43465
$ws.Cells.Item(7, 5).Value2 = 43373
$ws.Cells.Item(8, 4).Value2 = 146500
$ws.Cells.Item(8, 5).Value2 = 142600
$ws.Cells.Item(9, 4).Value2 = 100700
$ws.Cells.Item(9, 5).Value2 = 94100
$ws.Cells.Item(10, 4).Value2 = 45800
$ws.Cells.Item(10, 5).Value2 = 48500
$ws.Cells.Item(12, 4).Value2 = "NA"
$ws.Cells.Item(12, 5).Value2 = "NA"
$ws.Cells.Item(13, 4).Value2 = 0
$ws.Cells.Item(13, 5).Value2 = 0
$ws.Cells.Item(14, 4).Value2 = 0
$ws.Cells.Item(14, 5).Value2 = 0
$ws.Cells.Item(15, 4).Value2 = 0
$ws.Cells.Item(15, 5).Value2 = 0
$ws.Cells.Item(17, 4).Value2 = 114500
$ws.Cells.Item(17, 5).Value2 = 105500
$ws.Cells.Item(18, 4).Value2 = 32000
$ws.Cells.Item(18, 5).Value2 = 37100
$ws.Cells.Item(20, 4).Value2 = -3000
$ws.Cells.Item(20, 5).Value2 = 600
$ws.Cells.Item(21, 4).Value2 = "NA"
$ws.Cells.Item(21, 5).Value2 = "NA"
$ws.Cells.Item(22, 4).Value2 = 0
$ws.Cells.Item(22, 5).Value2 = 0
$ws.Cells.Item(23, 4).Value2 = 29000
$ws.Cells.Item(23, 5).Value2 = 37700
$ws.Cells.Item(24, 4).Value2 = 6200
$ws.Cells.Item(24, 5).Value2 = 8000
$ws.Cells.Item(25, 4).Value2 = 0
$ws.Cells.Item(25, 5).Value2 = 0
$ws.Cells.Item(26, 4).Value2 = 22900
$ws.Cells.Item(26, 5).Value2 = 29800
$ws.Cells.Item(27, 4).Value2 = 19600
$ws.Cells.Item(27, 5).Value2 = 26500
$ws.Cells.Item(28, 4).Value2 = 0
$ws.Cells.Item(28, 5).Value2 = 0
$ws.Cells.Item(29, 4).Value2 = 0
$ws.Cells.Item(29, 5).Value2 = "NA"
$ws.Cells.Item(30, 4).Value2 = 0
$ws.Cells.Item(30, 5).Value2 = 0
$ws.Cells.Item(31, 4).Value2 = 0
$ws.Cells.Item(31, 5).Value2 = 0
$ws.Cells.Item(32, 4).Value2 = 3000
$ws.Cells.Item(32, 5).Value2 = -600
$ws.Cells.Item(33, 4).Value2 = 19600
$ws.Cells.Item(33, 5).Value2 = 26500
$ws.Cells.Item(34, 4).Value2 = 0
$ws.Cells.Item(34, 5).Value2 = 0
$ws.Cells.Item(35, 4).Value2 = 19600
$ws.Cells.Item(35, 5).Value2 = 26500
$ws.Cells.Item(38, 4).Value2 = 43465
$ws.Cells.Item(38, 5).Value2 = 43373
$ws.Cells.Item(41, 4).Value2 = 425300
$ws.Cells.Item(41, 5).Value2 = 436200
$ws.Cells.Item(42, 4).Value2 = 0
$ws.Cells.Item(42, 5).Value2 = 0
$ws.Cells.Item(43, 4).Value2 = 220400
$ws.Cells.Item(43, 5).Value2 = 175900
$ws.Cells.Item(44, 4).Value2 = 0
$ws.Cells.Item(44, 5).Value2 = 0
$ws.Cells.Item(45, 4).Value2 = 9400
$ws.Cells.Item(45, 5).Value2 = 45200
$ws.Cells.Item(46, 4).Value2 = 0
$ws.Cells.Item(46, 5).Value2 = 0
$ws.Cells.Item(47, 4).Value2 = 14020300
$ws.Cells.Item(47, 5).Value2 = 13924000
$ws.Cells.Item(48, 4).Value2 = 0
$ws.Cells.Item(48, 5).Value2 = 0
$ws.Cells.Item(49, 4).Value2 = 0
$ws.Cells.Item(49, 5).Value2 = 0
$ws.Cells.Item(50, 4).Value2 = 0
$ws.Cells.Item(50, 5).Value2 = 0
$ws.Cells.Item(51, 4).Value2 = 0
$ws.Cells.Item(51, 5).Value2 = 0
$ws.Cells.Item(52, 4).Value2 = 6400
$ws.Cells.Item(52, 5).Value2 = 0
$ws.Cells.Item(53, 4).Value2 = 0
$ws.Cells.Item(53, 5).Value2 = 0
$ws.Cells.Item(54, 4).Value2 = 18694300
$ws.Cells.Item(54, 5).Value2 = 18474000
$ws.Cells.Item(57, 4).Value2 = 11900
$ws.Cells.Item(57, 5).Value2 = 260800
$ws.Cells.Item(58, 4).Value2 = 7757100
$ws.Cells.Item(58, 5).Value2 = 7378900
$ws.Cells.Item(59, 4).Value2 = 96700
$ws.Cells.Item(59, 5).Value2 = 87400
$ws.Cells.Item(60, 4).Value2 = 0
$ws.Cells.Item(60, 5).Value2 = 0
$ws.Cells.Item(61, 4).Value2 = 10015600
$ws.Cells.Item(61, 5).Value2 = 9906200
$ws.Cells.Item(62, 4).Value2 = 0
$ws.Cells.Item(62, 5).Value2 = 4600
$ws.Cells.Item(63, 4).Value2 = 0
$ws.Cells.Item(63, 5).Value2 = 0
$ws.Cells.Item(64, 4).Value2 = 0
$ws.Cells.Item(64, 5).Value2 = 0
$ws.Cells.Item(65, 4).Value2 = 0
$ws.Cells.Item(65, 5).Value2 = 0
$ws.Cells.Item(66, 4).Value2 = 17941800
$ws.Cells.Item(66, 5).Value2 = 17696400
$ws.Cells.Item(68, 4).Value2 = 0
$ws.Cells.Item(68, 5).Value2 = 0
$ws.Cells.Item(69, 4).Value2 = 0
$ws.Cells.Item(69, 5).Value2 = 0
$ws.Cells.Item(70, 4).Value2 = 204800
$ws.Cells.Item(70, 5).Value2 = 204800
$ws.Cells.Item(71, 4).Value2 = 0
$ws.Cells.Item(71, 5).Value2 = 0
$ws.Cells.Item(72, 4).Value2 = 393400
$ws.Cells.Item(72, 5).Value2 = 380000
$ws.Cells.Item(73, 4).Value2 = 0
$ws.Cells.Item(73, 5).Value2 = 0
$ws.Cells.Item(74, 4).Value2 = 0
$ws.Cells.Item(74, 5).Value2 = 0
$ws.Cells.Item(75, 4).Value2 = 0
$ws.Cells.Item(75, 5).Value2 = 0
$ws.Cells.Item(76, 4).Value2 = 547800
$ws.Cells.Item(76, 5).Value2 = 572800
$ws.Cells.Item(77, 4).Value2 = 0
$ws.Cells.Item(77, 5).Value2 = 0
$ws.Cells.Item(80, 4).Value2 = 43465
$ws.Cells.Item(80, 5).Value2 = 43373
$ws.Cells.Item(81, 4).Value2 = 19600
$ws.Cells.Item(81, 5).Value2 = 26500
$ws.Cells.Item(83, 4).Value2 = 0
$ws.Cells.Item(83, 5).Value2 = 0
$ws.Cells.Item(84, 4).Value2 = 0
$ws.Cells.Item(84, 5).Value2 = 0
$ws.Cells.Item(85, 4).Value2 = 0
$ws.Cells.Item(85, 5).Value2 = 0
$ws.Cells.Item(86, 4).Value2 = 0
$ws.Cells.Item(86, 5).Value2 = 0
$ws.Cells.Item(87, 4).Value2 = 0
$ws.Cells.Item(87, 5).Value2 = 0
$ws.Cells.Item(88, 4).Value2 = 0
$ws.Cells.Item(88, 5).Value2 = 0
$ws.Cells.Item(89, 4).Value2 = -52200
$ws.Cells.Item(89, 5).Value2 = 90100
$ws.Cells.Item(91, 4).Value2 = 0
$ws.Cells.Item(91, 5).Value2 = 0
$ws.Cells.Item(92, 4).Value2 = 0
$ws.Cells.Item(92, 5).Value2 = 0
$ws.Cells.Item(93, 4).Value2 = 0
$ws.Cells.Item(93, 5).Value2 = 0
$ws.Cells.Item(94, 4).Value2 = -359400
$ws.Cells.Item(94, 5).Value2 = 369300
$ws.Cells.Item(96, 4).Value2 = -9500
$ws.Cells.Item(96, 5).Value2 = -9500
$ws.Cells.Item(97, 4).Value2 = 0
$ws.Cells.Item(97, 5).Value2 = 0
$ws.Cells.Item(98, 4).Value2 = 0
$ws.Cells.Item(98, 5).Value2 = 0
$ws.Cells.Item(99, 4).Value2 = 0
$ws.Cells.Item(99, 5).Value2 = 0
$ws.Cells.Item(100, 4).Value2 = 400700
$ws.Cells.Item(100, 5).Value2 = -454000
$ws.Cells.Item(101, 4).Value2 = 0
$ws.Cells.Item(101, 5).Value2 = 0
$ws.Cells.Item(102, 4).Value2 = -10900
$ws.Cells.Item(102, 5).Value2 = 5300
